# Rebuild the "Estado de Cuenta" worker mora table (rows 16-47).
# Matches commit: "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# - YORBIS (73508103): row 16 Valor Mora changes 79075 -> 47445
# - rows 25-32 duplicate YORBIS periods 2406..2311 (reverse chronological)
# - rows 33-39 become GERMAN MARTINEZ CALDERON (19596360), periods 2503..2409
# - rows 40-47 become GENDRIS ZULEIMA OROZCO RODRIGUEZ (1065625378), periods 2503..2408

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "73508103"
$ws.Range("D16").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E16").Value = "2311"
$ws.Range("F16").Value = 47445
$ws.Range("G16").Value = 1976894
$ws.Range("C17").Value = "73508103"
$ws.Range("D17").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E17").Value = "2312"
$ws.Range("F17").Value = 79075
$ws.Range("G17").Value = 1976894
$ws.Range("C18").Value = "73508103"
$ws.Range("D18").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E18").Value = "2401"
$ws.Range("F18").Value = 79075
$ws.Range("G18").Value = 1976894
$ws.Range("C19").Value = "73508103"
$ws.Range("D19").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E19").Value = "2402"
$ws.Range("F19").Value = 79075
$ws.Range("G19").Value = 1976894
$ws.Range("C20").Value = "73508103"
$ws.Range("D20").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E20").Value = "2403"
$ws.Range("F20").Value = 79075
$ws.Range("G20").Value = 1976894
$ws.Range("C21").Value = "73508103"
$ws.Range("D21").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E21").Value = "2404"
$ws.Range("F21").Value = 79075
$ws.Range("G21").Value = 1976894
$ws.Range("C22").Value = "73508103"
$ws.Range("D22").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E22").Value = "2405"
$ws.Range("F22").Value = 79075
$ws.Range("G22").Value = 1976894
$ws.Range("C23").Value = "73508103"
$ws.Range("D23").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E23").Value = "2406"
$ws.Range("F23").Value = 79075
$ws.Range("G23").Value = 1976894
$ws.Range("C24").Value = "73508103"
$ws.Range("D24").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E24").Value = "2407"
$ws.Range("F24").Value = 79075
$ws.Range("G24").Value = 1976894
$ws.Range("C25").Value = "73508103"
$ws.Range("D25").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E25").Value = "2406"
$ws.Range("F25").Value = 79075
$ws.Range("G25").Value = 1976894
$ws.Range("C26").Value = "73508103"
$ws.Range("D26").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E26").Value = "2405"
$ws.Range("F26").Value = 79075
$ws.Range("G26").Value = 1976894
$ws.Range("C27").Value = "73508103"
$ws.Range("D27").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E27").Value = "2404"
$ws.Range("F27").Value = 79075
$ws.Range("G27").Value = 1976894
$ws.Range("C28").Value = "73508103"
$ws.Range("D28").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E28").Value = "2403"
$ws.Range("F28").Value = 79075
$ws.Range("G28").Value = 1976894
$ws.Range("C29").Value = "73508103"
$ws.Range("D29").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E29").Value = "2402"
$ws.Range("F29").Value = 79075
$ws.Range("G29").Value = 1976894
$ws.Range("C30").Value = "73508103"
$ws.Range("D30").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E30").Value = "2401"
$ws.Range("F30").Value = 79075
$ws.Range("G30").Value = 1976894
$ws.Range("C31").Value = "73508103"
$ws.Range("D31").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E31").Value = "2312"
$ws.Range("F31").Value = 79075
$ws.Range("G31").Value = 1976894
$ws.Range("C32").Value = "73508103"
$ws.Range("D32").Value = "YORBIS ANTONIO ROSADO MENDOZA"
$ws.Range("E32").Value = "2311"
$ws.Range("F32").Value = 79075
$ws.Range("G32").Value = 1976894
$ws.Range("C33").Value = "19596360"
$ws.Range("D33").Value = "GERMAN MARTINEZ CALDERON"
$ws.Range("E33").Value = "2503"
$ws.Range("F33").Value = 108000
$ws.Range("G33").Value = 4500000
$ws.Range("C34").Value = "19596360"
$ws.Range("D34").Value = "GERMAN MARTINEZ CALDERON"
$ws.Range("E34").Value = "2502"
$ws.Range("F34").Value = 180000
$ws.Range("G34").Value = 4500000
$ws.Range("C35").Value = "19596360"
$ws.Range("D35").Value = "GERMAN MARTINEZ CALDERON"
$ws.Range("E35").Value = "2501"
$ws.Range("F35").Value = 180000
$ws.Range("G35").Value = 4500000
$ws.Range("C36").Value = "19596360"
$ws.Range("D36").Value = "GERMAN MARTINEZ CALDERON"
$ws.Range("E36").Value = "2412"
$ws.Range("F36").Value = 180000
$ws.Range("G36").Value = 4500000
$ws.Range("C37").Value = "19596360"
$ws.Range("D37").Value = "GERMAN MARTINEZ CALDERON"
$ws.Range("E37").Value = "2411"
$ws.Range("F37").Value = 180000
$ws.Range("G37").Value = 4500000
$ws.Range("C38").Value = "19596360"
$ws.Range("D38").Value = "GERMAN MARTINEZ CALDERON"
$ws.Range("E38").Value = "2410"
$ws.Range("F38").Value = 180000
$ws.Range("G38").Value = 4500000
$ws.Range("C39").Value = "19596360"
$ws.Range("D39").Value = "GERMAN MARTINEZ CALDERON"
$ws.Range("E39").Value = "2409"
$ws.Range("F39").Value = 180000
$ws.Range("G39").Value = 4500000
$ws.Range("C40").Value = "1065625378"
$ws.Range("D40").Value = "GENDRIS ZULEIMA OROZCO RODRIGUEZ"
$ws.Range("E40").Value = "2503"
$ws.Range("F40").Value = 72000
$ws.Range("G40").Value = 3000000
$ws.Range("C41").Value = "1065625378"
$ws.Range("D41").Value = "GENDRIS ZULEIMA OROZCO RODRIGUEZ"
$ws.Range("E41").Value = "2502"
$ws.Range("F41").Value = 120000
$ws.Range("G41").Value = 3000000
$ws.Range("C42").Value = "1065625378"
$ws.Range("D42").Value = "GENDRIS ZULEIMA OROZCO RODRIGUEZ"
$ws.Range("E42").Value = "2501"
$ws.Range("F42").Value = 120000
$ws.Range("G42").Value = 3000000
$ws.Range("C43").Value = "1065625378"
$ws.Range("D43").Value = "GENDRIS ZULEIMA OROZCO RODRIGUEZ"
$ws.Range("E43").Value = "2412"
$ws.Range("F43").Value = 120000
$ws.Range("G43").Value = 3000000
$ws.Range("C44").Value = "1065625378"
$ws.Range("D44").Value = "GENDRIS ZULEIMA OROZCO RODRIGUEZ"
$ws.Range("E44").Value = "2411"
$ws.Range("F44").Value = 120000
$ws.Range("G44").Value = 3000000
$ws.Range("C45").Value = "1065625378"
$ws.Range("D45").Value = "GENDRIS ZULEIMA OROZCO RODRIGUEZ"
$ws.Range("E45").Value = "2410"
$ws.Range("F45").Value = 120000
$ws.Range("G45").Value = 3000000
$ws.Range("C46").Value = "1065625378"
$ws.Range("D46").Value = "GENDRIS ZULEIMA OROZCO RODRIGUEZ"
$ws.Range("E46").Value = "2409"
$ws.Range("F46").Value = 120000
$ws.Range("G46").Value = 3000000
$ws.Range("C47").Value = "1065625378"
$ws.Range("D47").Value = "GENDRIS ZULEIMA OROZCO RODRIGUEZ"
$ws.Range("E47").Value = "2408"
$ws.Range("F47").Value = 120000
$ws.Range("G47").Value = 3000000
